$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Title "
$ws.Range("B1").Value = "Studio "
$ws.Range("C1").Value = "Gen "
$ws.Range("D1").Value = "Likes"

$ws.Range("A2").Value = "Avatar"
$ws.Range("B2").Value = "Mara"
$ws.Range("C2").Value = "Fantasy"
$ws.Range("D2").Value = 123

$ws.Range("A3").Value = "Avengers"
$ws.Range("B3").Value = "Marvel"
$ws.Range("C3").Value = "Fantasy"
$ws.Range("D3").Value = 4998

$ws.Range("A4").Value = "Wednesday"
$ws.Range("B4").Value = "Roman"
$ws.Range("C4").Value = "Drama"
$ws.Range("D4").Value = 234

$ws.Range("D4").Select()
